# "Added CDS All studies testcase"
#
# The SamplesTab query (row 3, column B on Sheet1) is rewritten: the old
# "Sample ID" query (which pulled Tumor / Analyte Type columns in addition
# to Sample ID / Participant ID / Study Name / Accession) is replaced with a
# leaner version that only selects Sample ID / Participant ID / Study Name /
# Accession. The ParticipantsTab (B2) and FilesTab (B4) queries are
# unchanged in content.
#
# The sheet's view is also nudged: the window is scrolled down one row and
# the active selection moves from A2 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
  s.phs_accession = 'phs001437' AND gi.library_strategy = 'WXS'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# Replace the SamplesTab (row 3) TabQuery cell with the new, narrower query.
$ws.Range("B3").Value = $newSamplesQuery

# Update the view: scroll the window down a row and move the selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C3").Select()
